$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header
$ws.Range("D1").Value = "type_order"

# Add value to existing row 2
$ws.Range("D2").Value = "SALE"

# Add new row 3 with data
$ws.Range("A3").Value = "DELETED"
$ws.Range("B3").Value = "now()"
$ws.Range("C3").Value = 1212121
$ws.Range("D3").Value = "SALE"

# Update selection to match target
$ws.Range("E7").Select()
